$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs target) - updated TPM-derived values
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ntn4"
$ws.Cells.Item(2, 3).Value = "Dcc"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.497392666666667
$ws.Cells.Item(2, 8).Value = 4.492178
$ws.Cells.Item(2, 9).Value = 0.07263102411830044
$ws.Cells.Item(2, 10).Value = 0.07263102411830046
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.013316
$ws.Cells.Item(2, 14).Value = 0.039948
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(2, 16).Value = 1
$ws.Cells.Item(2, 17).Value = 0.01993928074933333
$ws.Cells.Item(2, 18).Value = 0.179453526744
$ws.Cells.Item(2, 19).Value = 0.07263102411830044
$ws.Cells.Item(2, 20).Value = 0.07263102411830046

# Row 3 (FAPs target) - updated TPM-derived values
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Ntn4"
$ws.Cells.Item(3, 3).Value = "Dcc"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 7.438224333333333
$ws.Cells.Item(3, 8).Value = 22.314673
$ws.Cells.Item(3, 9).Value = 0.3607910356301526
$ws.Cells.Item(3, 10).Value = 0.3607910356301526
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.013316
$ws.Cells.Item(3, 14).Value = 0.039948
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(3, 16).Value = 1
$ws.Cells.Item(3, 17).Value = 0.09904739522266666
$ws.Cells.Item(3, 18).Value = 0.8914265570039999
$ws.Cells.Item(3, 19).Value = 0.3607910356301526
$ws.Cells.Item(3, 20).Value = 0.3607910356301526

# Row 4 (Inflammatory-Mac target, newly inserted cluster label) - updated TPM-derived values
$ws.Cells.Item(4, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 2).Value = "Ntn4"
$ws.Cells.Item(4, 3).Value = "Dcc"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1679733333333333
$ws.Cells.Item(4, 8).Value = 0.5039199999999999
$ws.Cells.Item(4, 9).Value = 0.008147545728084229
$ws.Cells.Item(4, 10).Value = 0.008147545728084229
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.013316
$ws.Cells.Item(4, 14).Value = 0.039948
$ws.Cells.Item(4, 15).Value = 1
$ws.Cells.Item(4, 16).Value = 1
$ws.Cells.Item(4, 17).Value = 0.002236732906666666
$ws.Cells.Item(4, 18).Value = 0.02013059616
$ws.Cells.Item(4, 19).Value = 0.008147545728084229
$ws.Cells.Item(4, 20).Value = 0.008147545728084229

# Row 5 (MuSCs target) - updated TPM-derived values
$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Ntn4"
$ws.Cells.Item(5, 3).Value = "Dcc"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 11.470741
$ws.Cells.Item(5, 8).Value = 34.412223
$ws.Cells.Item(5, 9).Value = 0.5563882372152958
$ws.Cells.Item(5, 10).Value = 0.5563882372152958
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.013316
$ws.Cells.Item(5, 14).Value = 0.039948
$ws.Cells.Item(5, 15).Value = 1
$ws.Cells.Item(5, 16).Value = 1
$ws.Cells.Item(5, 17).Value = 0.152744387156
$ws.Cells.Item(5, 18).Value = 1.374699484404
$ws.Cells.Item(5, 19).Value = 0.5563882372152958
$ws.Cells.Item(5, 20).Value = 0.5563882372152958

# Row 6 (Resolving-Mac target) - new row added at the end of the table
$ws.Cells.Item(6, 1).Value = "Resolving-Mac"
$ws.Cells.Item(6, 2).Value = "Ntn4"
$ws.Cells.Item(6, 3).Value = "Dcc"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.042102
$ws.Cells.Item(6, 8).Value = 0.126306
$ws.Cells.Item(6, 9).Value = 0.002042157308166786
$ws.Cells.Item(6, 10).Value = 0.002042157308166786
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.013316
$ws.Cells.Item(6, 14).Value = 0.039948
$ws.Cells.Item(6, 15).Value = 1
$ws.Cells.Item(6, 16).Value = 1
$ws.Cells.Item(6, 17).Value = 0.0005606302319999999
$ws.Cells.Item(6, 18).Value = 0.005045672088
$ws.Cells.Item(6, 19).Value = 0.002042157308166786
$ws.Cells.Item(6, 20).Value = 0.002042157308166786
